# Apply the "Part 1c -> Part 1cb3" update described in the commit:
#   updated part1c, now use part1cb3, part1c to be removed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Part 1")

# ---------------------------------------------------------------------------
# 1. Update the Solver (Location A turbine counts) model so it now targets
#    L21 (Total Investment) instead of the Location B turbine row, and the
#    investment cap moves from a literal "integer" placeholder to $K$6.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Part 1!solver_lhs5") {
        $n.RefersTo = "='Part 1'!`$L`$21"
    }
    elseif ($n.Name -eq "Part 1!solver_num") {
        $n.RefersTo = "=5"
    }
    elseif ($n.Name -eq "Part 1!solver_rel5") {
        $n.RefersTo = "=1"
    }
    elseif ($n.Name -eq "Part 1!solver_rhs5") {
        $n.RefersTo = "='Part 1'!`$K`$6"
    }
}

# ---------------------------------------------------------------------------
# 2. Updated Location A turbine mix (F18/G18)
# ---------------------------------------------------------------------------
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 11

# ---------------------------------------------------------------------------
# 3. New "Cost of Investment" columns (K/L), rows 18-21
# ---------------------------------------------------------------------------
$ws.Range("K18").Value = "Cost of Investment Location A"
$ws.Range("L18").Formula = "=SUMPRODUCT(B10:D10,Location_A_turbines)"

$ws.Range("K19").Value = "Cost of Investment Location B"
$ws.Range("L19").Formula = "=SUMPRODUCT(B11:D11,Location_A_turbines)"

$ws.Range("K20").Value = "Cost of Investment Location C"
$ws.Range("L20").Formula = "=SUMPRODUCT(B12:D12,Location_A_turbines)"

$ws.Range("K21").Value = "Total Investment"
$ws.Range("L21").Formula = "=SUM(L18:L20)"

# Match formatting of the analogous existing cells.
$ws.Range("J2").Copy()
$ws.Range("K18:K20").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("J21").Copy()
$ws.Range("L18:L21").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("K5").Copy()
$ws.Range("K21").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("K21").Borders.LineStyle = -4142 # xlLineStyleNone

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Return on Investment formulas now also subtract a share of the new
#    total investment cap ($K$6/3)
# ---------------------------------------------------------------------------
$ws.Range("F25").Formula = "=SUMPRODUCT(B7:D7,Location_A_turbines)-((SUMPRODUCT(B10:D10,Location_A_turbines)+(SUMPRODUCT(B13:D13,Location_A_turbines)))/((SUMPRODUCT(B10:D10,Location_A_turbines)+(SUMPRODUCT(B13:D13,Location_A_turbines))))) -(`$K`$6/3)"
$ws.Range("F26").Formula = "=SUMPRODUCT(B8:D8,Location_A_turbines)-((SUMPRODUCT(B11:D11,Location_A_turbines)+(SUMPRODUCT(B14:D14,Location_A_turbines)))/((SUMPRODUCT(B11:D11,Location_A_turbines)+(SUMPRODUCT(B14:D14,Location_A_turbines))))) -(`$K`$6/3)"
$ws.Range("F27").Formula = "=SUMPRODUCT(B9:D9,Location_A_turbines)-((SUMPRODUCT(B12:D12,Location_A_turbines)+(SUMPRODUCT(B15:D15,Location_A_turbines)))/((SUMPRODUCT(B12:D12,Location_A_turbines)+(SUMPRODUCT(B15:D15,Location_A_turbines))))) -(`$K`$6/3)"

# Leftover empty formatted cells on rows 25/26 column C (picked up by the
# author while reformatting the ROI block) - same currency format as the
# blank K7:K9 placeholder cells.
$ws.Range("K7").Copy()
$ws.Range("C25:C26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()

# ---------------------------------------------------------------------------
# 5. Column width tweaks to fit the new content
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 55.5703125
$ws.Columns.Item(11).ColumnWidth = 31.7109375
$ws.Columns.Item(12).ColumnWidth = 15

# ---------------------------------------------------------------------------
# 6. Remove the embedded picture that used to sit over the ROI table
# ---------------------------------------------------------------------------
while ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.Item(1).Delete()
}

# ---------------------------------------------------------------------------
# 7. Sheet view: zoomed to 85%, selection moved to G29
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 85
$ws.Range("G29").Select() | Out-Null

Write-Host "Edit complete"
